# Updates the "展览" (Exhibitions) and "全部类型" (All Types) sheets of the
# Jiangxi comic-convention tracker workbook to the data snapshot generated
# at commit 456a3b4:
#   - bump several "want to go" counters (column F) for existing rows
#   - insert a brand-new row for 上饶·IX Group as the new row 32 (shifting
#     everything below it down by one)
#   - append a brand-new row 38 for 南昌·第四届龙年动漫展 at the bottom
#
# Both worksheets hold identical tables, so the same sequence of edits is
# replayed on each of them.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # ---- 1. Straightforward "want to go" (column F) counter bumps ----
    $ws.Range("F7").Value  = 146
    $ws.Range("F9").Value  = 22
    $ws.Range("F13").Value = 26
    $ws.Range("F15").Value = 553
    $ws.Range("F17").Value = 480
    $ws.Range("F18").Value = 412
    $ws.Range("F23").Value = 1198
    $ws.Range("F24").Value = 2873
    $ws.Range("F27").Value = 549
    $ws.Range("F28").Value = 55
    $ws.Range("F29").Value = 1619
    $ws.Range("F30").Value = 565

    # ---- 2. Insert the new row 32 (上饶·IX Group) ----
    # This shifts the old rows 32-36 down to 33-37, carrying their values,
    # so only the counters that actually changed (F33, F34, F36) need a
    # follow-up edit below.
    $ws.Rows("32:32").Insert()

    # Column A in this table is a plain hand-typed row index (not a
    # formula), styled bold/centered/bordered (style of A31/A33/...).
    # Copy that formatting onto the freshly inserted, still-blank A32 so
    # the new row matches its neighbours before we fill in the value.
    $ws.Range("A31").Copy()
    $ws.Range("A32").PasteSpecial(-4122)

    $ws.Range("A32").Value = 31

    # Dates typed as plain "YYYY-MM-DD" text get auto-parsed as real dates
    # by Excel's normal Value assignment; format the cell as Text first so
    # it is stored as the literal string, matching the rest of column B.
    $ws.Range("B32").NumberFormat = "@"
    $ws.Range("B32").Value = "2024-08-04"

    $ws.Range("C32").Value = "上饶·第十五届IX Group国风嘉年华暨十周年庆典"
    $ws.Range("D32").Value = "高铁经济试验区凤凰东大道16号7幢 上饶饶商金茂诚悦酒店(上饶高铁站)"
    $ws.Range("E32").Value = "2024.08.04 09:30-08.04 17:30"
    $ws.Range("F32").Value = 6
    $ws.Range("G32").Value = 60
    $ws.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=87225"
    $ws.Range("I32").Value = "//i2.hdslb.com/bfs/openplatform/202406/l5fIXZSX1717562269098.jpeg"

    # ---- 3. Counter bumps on the rows that shifted down ----
    $ws.Range("F33").Value = 268   # 九江·第一届异次元动漫嘉年华 (was F32=265)
    $ws.Range("F34").Value = 394   # 南昌·第一届异次元动漫嘉年华 (was F33=391)
    $ws.Range("F36").Value = 607   # 赣州·第二届异次元动漫嘉年华 (was F35=605)

    # ---- 4. Append the new row 38 (南昌·第四届龙年动漫展) ----
    $ws.Range("A37").Copy()
    $ws.Range("A38").PasteSpecial(-4122)

    $ws.Range("A38").Value = 37

    $ws.Range("B38").NumberFormat = "@"
    $ws.Range("B38").Value = "2024-08-24"

    $ws.Range("C38").Value = "南昌·第四届龙年动漫展——暑假最后的狂欢"
    $ws.Range("D38").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
    $ws.Range("E38").Value = "2024.08.24 10:00-08.25 18:00"
    $ws.Range("F38").Value = 1
    $ws.Range("G38").Value = 45
    $ws.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=87135"
    $ws.Range("I38").Value = "//i0.hdslb.com/bfs/openplatform/202406/mDtqZeQd1718033555304.jpeg"
}
